$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header I2: "QF" -> "B field"
$ws.Range("I2").Value = "B field"

# Add new header J2: "E field"
$ws.Range("J2").Value = "E field"

# Set column J width to match column I
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# Fill J3:J13 with the new formula (relative reference on I column row)
$ws.Range("J3").Formula = "=(`$F`$20/(I3^2-`$F`$20))-(`$F`$16/(I3^2-`$F`$16))"
$ws.Range("J4:J13").Formula = "=(`$F`$20/(I4^2-`$F`$20))-(`$F`$16/(I4^2-`$F`$16))"

# Apply scientific number format to O7
$ws.Range("O7").NumberFormat = "0.00E+00"

# Restore selection to match target state (closest achievable: range H1:J13)
$ws.Range("H1:J13").Select() | Out-Null
